# Updates the cryptos price (D) and 1h-volume-change (E) columns with
# refreshed values, as produced by the "Updated cryptos list ... with
# GitHub Actions" commit. All of these cells hold plain text (price
# strings use "." as both thousands- and decimal-separators, e.g.
# "67.674.33", and the % cells are padded with spaces), so values are
# written as text rather than being left for Excel to auto-infer as
# numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text does not look like a plain number to Excel's
# auto-detection (multi-dot prices, subscript-digit prices, and every
# "  +/-X.XX%  " cell) can be written directly.
$plainTextUpdates = @{
    "D2" = '67.674.33'
    "E2" = '  -0.24%  '
    "D3" = '2.425.89'
    "E3" = '  +0.41%  '
    "E4" = '  +0.04%  '
    "E6" = '  +0.41%  '
    "E7" = '  +0.01%  '
    "E8" = '  +2.92%  '
    "E9" = '  +8.31%  '
    "E10" = '  -0.03%  '
    "E11" = '  -0.38%  '
    "E12" = '  +1.26%  '
    "D13" = '67.623.94'
    "E13" = '  -0.12%  '
    "E14" = '  +2.26%  '
    "E15" = '  +0.06%  '
    "E16" = '  -2.10%  '
    "E17" = '  -0.46%  '
    "E18" = '  -0.61%  '
    "E19" = '  +1.30%  '
    "E20" = '  +0.06%  '
    "E21" = '  +2.05%  '
    "E22" = '  +0.72%  '
    "E23" = '  +1.15%  '
    "D25" = '0.0₃0805'
    "E25" = '  +1.80%  '
    "E26" = '  +0.94%  '
    "E27" = '  +0.13%  '
    "E28" = '  -1.82%  '
    "E29" = '  +2.84%  '
    "E30" = '  +0.54%  '
    "E31" = '  +2.53%  '
    "E32" = '  -0.27%  '
    "E33" = '  +0.12%  '
    "E34" = '  +1.29%  '
    "E35" = '  -3.24%  '
    "E36" = '  -1.09%  '
    "E37" = '  -1.57%  '
    "E38" = '  +2.57%  '
    "E39" = '  -0.16%  '
    "E40" = '  +1.04%  '
    "E41" = '  +1.61%  '
    "E42" = '  -0.55%  '
    "E43" = '  +0.36%  '
    "E44" = '  +0.92%  '
    "E45" = '  +0.96%  '
    "E46" = '  +1.75%  '
    "E47" = '  +0.73%  '
    "E48" = '  -4.86%  '
    "E49" = '  +0.68%  '
    "D50" = '0.0₆0203'
    "E50" = '  +5.47%  '
    "E51" = '  +1.83%  '
}

# Cells whose new text DOES look like an ordinary number (e.g. "552.50",
# "0.509") need a leading apostrophe so Excel stores them as text -
# exactly as it would if a person typed an apostrophe-prefixed entry -
# instead of silently converting them to numeric values and dropping
# formatting such as trailing zeros.
$forcedTextUpdates = @{
    "D5" = '552.50'
    "D6" = '159.41'
    "D8" = '0.509'
    "D11" = '0.327'
    "D12" = '4.78'
    "D15" = '22.95'
    "D16" = '10.32'
    "D17" = '333.38'
    "D18" = '6.81'
    "D19" = '3.77'
    "D22" = '66.14'
    "D23" = '3.61'
    "D24" = '8.06'
    "D26" = '7.07'
    "D28" = '418.85'
    "D31" = '160.85'
    "D32" = '18.92'
    "D34" = '17.77'
    "D35" = '0.104'
    "D37" = '4.24'
    "D39" = '1.06'
    "D40" = '1.99'
    "D41" = '3.32'
    "D42" = '128.51'
    "D43" = '0.0708'
    "D44" = '0.477'
    "D45" = '0.553'
    "D46" = '0.0910'
    "D48" = '1.33'
    "D49" = '16.54'
    "D51" = '0.0427'
}

foreach ($addr in $plainTextUpdates.Keys) {
    $ws.Range($addr).Value = $plainTextUpdates[$addr]
}

foreach ($addr in $forcedTextUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $forcedTextUpdates[$addr]
}
